$d = $word.ActiveDocument

# Manual line breaks (<w:br/>) are represented as vertical-tab (chr 11)
# characters in Word's Range.Text / Find & Replace.
$vt = [char]11
$apos = [char]0x2019

# 1) Update the "Reviewed on" date.
$d.Content.Find.Execute(
    "Reviewed on: 2022-12-08", $true, $false, $false, $false, $false,
    $true, 1, $false, "Reviewed on: 2022-12-09", 2) | Out-Null

# 2) Add one extra manual line break (total 4 instead of 3) after each of
#    the ten English example sentences.
$sentences = @(
    "he purposely put sand in other people" + $apos + "s food.",
    "she listens to what I have to say.",
    "until they realize how dangerous she can be.",
    "but somehow it just tasted rotten.",
    "crocodiles into the ocean.",
    "that my baby brother was a jewel thief.",
    "has made creating trash more acceptable.",
    "with his hazard lights on.",
    "and wore it on a sunny day.",
    "and screams for more memes."
)

foreach ($s in $sentences) {
    $old = $s + $vt + $vt + $vt
    $new = $s + $vt + $vt + $vt + $vt
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
        $true, 1, $false, $new, 2) | Out-Null
}
